# Update cryptocurrency price/volume data per scheduled GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.758.68'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '3.403.12'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '408.48'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.54%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '127.77'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -3.78%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.631'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +6.52%  '
$ws.Range('E8').Value = '  -0.38%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.726'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +5.92%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.138'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +9.08%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '42.26'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.37%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '9.04'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +7.12%  '
$ws.Range('D14').Value = '3.943.21'
$ws.Range('E14').Value = '  -0.78%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '21.18'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +6.97%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.0000201'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +41.54%  '
$ws.Range('D17').Value = '3.406.13'
$ws.Range('E17').Value = '  -0.61%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '12.02'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +4.67%  '
$ws.Range('E19').Value = '  +5.08%  '
$ws.Range('D20').Value = '61.728.95'
$ws.Range('E20').Value = '  -0.45%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '451.61'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +45.46%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '91.27'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +9.32%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '3.14'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.72%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '12.85'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.34%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '3.24'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +3.10%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '33.27'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +12.27%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '8.64'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +6.64%  '
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('E29').Value = '  -1.05%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '7.52'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -2.32%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '11.90'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +4.76%  '
$ws.Range('E32').Value = '  -3.40%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '42.59'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('E34').Value = '  -1.52%  '
$ws.Range('E35').Value = '  +0.11%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.0495'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +2.17%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '53.09'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('E40').Value = '  +6.99%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '2.90'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -0.76%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.313'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -4.11%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '140.93'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.01%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '4.19'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +5.78%  '
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('E46').Value = '  +8.14%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '16.42'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -1.73%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '22.30'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +5.14%  '
$ws.Range('D49').Value = '3.751.20'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').Value = '2.101.30'
$ws.Range('E50').Value = '  -0.74%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '105.89'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +27.74%  '
